$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 33, shifting existing rows 33-127 down to 34-128
$ws.Rows("33:33").Insert()

# Populate the newly inserted row 33 with the new data record
$ws.Range("A33").Value = 11
$ws.Range("B33").Value = "Vega Monumental Concepción"
$ws.Range("C33").Value = "Bíobío"
$ws.Range("D33").Value = 45099
$ws.Range("E33").Value = 8
$ws.Range("F33").Value = 100112012
$ws.Range("G33").Value = "Espinaca"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 220
$ws.Range("K33").Value = 6500
$ws.Range("L33").Value = 7000
$ws.Range("M33").Value = 6727
$ws.Range("N33").Value = "`$/cuna 10 kilos"
$ws.Range("O33").Value = "Región Metropolitana"
$ws.Range("P33").Value = 673
$ws.Range("Q33").Value = 10
$ws.Range("R33").Value = "Hortaliza"
